$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2024-05-25)
$ws.Range("B2").Value = 0.003994804209775715
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 0.6610771962923778

# Row 3 (2024-03-04)
$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 5.014808316549482
